$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells G1 ("cc_1st") and H1 ("cc_2nd"), centered horizontally & vertically.
$ws.Range("G1").Value = "cc_1st"
$ws.Range("H1").Value = "cc_2nd"
$ws.Range("G1:H1").HorizontalAlignment = -4108
$ws.Range("G1:H1").VerticalAlignment = -4108

# New cumulative case-count data in columns G (cc_1st) and H (cc_2nd).
$ws.Range("G6").Value = 226
$ws.Range("H6").Value = 113

$ws.Range("G7").VerticalAlignment = -4108

$ws.Range("G9").Value = 2670
$ws.Range("H9").Value = 2152

$ws.Range("G10").Value = 3054
$ws.Range("H10").Value = 3054

$ws.Range("G11").Value = 6542
$ws.Range("H11").Value = 5937

$ws.Range("G12").Value = 8982
$ws.Range("H12").Value = 6388

$ws.Range("G13").Value = 10728
$ws.Range("H13").Value = 7023

# Fill in the previously-blank daily counts for July 17 (row 14).
$ws.Range("C14").Value = 53
$ws.Range("D14").Value = 68
$ws.Range("E14").Value = 47
$ws.Range("F14").Value = 57

$ws.Range("G14").Value = 11914
$ws.Range("H14").Value = 7501

$ws.Range("G15").Value = 13708
$ws.Range("H15").Value = 8640

[void]$ws.Range("I17").Select()
